$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This single shared string is referenced from the Overview sheet
#    (columns zh-cn/de-de) as well as from the per-language "Status" column,
#    so a workbook-wide replace keeps every occurrence in sync.
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US") | Out-Null
}

# ---------------------------------------------------------------------------
# 2. Fill in the "handback" columns (Latest Target File / Latest Handback
#    File / Latest Handback DateTime) on the zh-cn and de-de sheets, and
#    link the new "Latest Target File" cells back to the source doc - same
#    as column A already does.
# ---------------------------------------------------------------------------
$a660Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a28874929815dbb84685f6797f48ed0b72518b6b/e2e/a660d210-72bf-4417-ba7b-500094e8d6fc.md"
$f865Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a28874929815dbb84685f6797f48ed0b72518b6b/e2e/f865af1c-6834-4d5e-ba98-fee8ccfcb67b.md"

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("I2").Value = "a660d210-72bf-4417-ba7b-500094e8d6fc.md"
$ws_zhcn.Range("J2").Value = "a660d210-72bf-4417-ba7b-500094e8d6fc.3fa8fd22c0907f245df1c39a257eb4b45218e267.zh-cn.xlf"
$ws_zhcn.Range("K2").Value = "2016-08-31 03:50:48"
$ws_zhcn.Range("I3").Value = "f865af1c-6834-4d5e-ba98-fee8ccfcb67b.md"
$ws_zhcn.Range("J3").Value = "f865af1c-6834-4d5e-ba98-fee8ccfcb67b.49d156d41494eefcc22441bba53874420bc60b1c.zh-cn.xlf"
$ws_zhcn.Range("K3").Value = "2016-08-31 03:50:48"

$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("I2"), $a660Url, [Type]::Missing, [Type]::Missing, "a660d210-72bf-4417-ba7b-500094e8d6fc.md") | Out-Null
$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("I3"), $f865Url, [Type]::Missing, [Type]::Missing, "f865af1c-6834-4d5e-ba98-fee8ccfcb67b.md") | Out-Null
$ws_zhcn.Range("I2").Style = "HyperLink"
$ws_zhcn.Range("I3").Style = "HyperLink"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("I2").Value = "a660d210-72bf-4417-ba7b-500094e8d6fc.md"
$ws_dede.Range("J2").Value = "a660d210-72bf-4417-ba7b-500094e8d6fc.3fa8fd22c0907f245df1c39a257eb4b45218e267.de-de.xlf"
$ws_dede.Range("K2").Value = "2016-08-31 03:51:06"
$ws_dede.Range("I3").Value = "f865af1c-6834-4d5e-ba98-fee8ccfcb67b.md"
$ws_dede.Range("J3").Value = "f865af1c-6834-4d5e-ba98-fee8ccfcb67b.49d156d41494eefcc22441bba53874420bc60b1c.de-de.xlf"
$ws_dede.Range("K3").Value = "2016-08-31 03:51:06"

$ws_dede.Hyperlinks.Add($ws_dede.Range("I2"), $a660Url, [Type]::Missing, [Type]::Missing, "a660d210-72bf-4417-ba7b-500094e8d6fc.md") | Out-Null
$ws_dede.Hyperlinks.Add($ws_dede.Range("I3"), $f865Url, [Type]::Missing, [Type]::Missing, "f865af1c-6834-4d5e-ba98-fee8ccfcb67b.md") | Out-Null
$ws_dede.Range("I2").Style = "HyperLink"
$ws_dede.Range("I3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# 3. Widen the columns that now hold longer text (report generated for
#    handback): Overview!E:F, and Status/Latest Target File/Latest
#    Handback File on both language sheets.
# ---------------------------------------------------------------------------
$ws_overview = $wb.Worksheets.Item("Overview")
$ws_overview.Columns.Item("E").ColumnWidth = 29.14437166849777
$ws_overview.Columns.Item("F").ColumnWidth = 29.14437166849777

$ws_zhcn.Columns.Item("C").ColumnWidth = 29.14437166849777
$ws_zhcn.Columns.Item("I").ColumnWidth = 39.166666666666664
$ws_zhcn.Columns.Item("J").ColumnWidth = 39.166666666666664

$ws_dede.Columns.Item("C").ColumnWidth = 29.14437166849777
$ws_dede.Columns.Item("I").ColumnWidth = 39.166666666666664
$ws_dede.Columns.Item("J").ColumnWidth = 39.166666666666664
